$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Caso de uso 3 (fila 4): renombrar otorgarTarifaSocial -> solicitarTarifaSocial ---
$ws.Range("C4").Value = "solicitarTarifaSocial"

# --- Caso de uso 4 (fila 5): solicitarBoletoEstudiantil ahora toma un parametro y devuelve void ---
$ws.Range("E5").Value = "boletoEstudiantil: BoletoEstudiantil"
$ws.Range("F5").Value = "void"

# --- Insertar una fila nueva antes de la fila 7 para separar MaquinaColectivo / MaquinaTren ---
$ws.Rows.Item(7).Insert()

# --- Caso de uso 5 (fila 6): pasa de Maquina.cobrar generico a MaquinaColectivo.cobrar ---
$ws.Range("B6").Value = "MaquinaColectivo"
$ws.Range("C6").Value = "cobrar"
$ws.Range("D6").Value = "descuenta el valor del boleto (Colectivo)"
$ws.Range("E6").Value = "tarjeta: Tarjeta, boletoColectivo: BoletoColectivo"
$ws.Range("F6").Value = "void"

# --- Caso de uso 6 (fila 7, nueva): MaquinaTren.cobrar ---
$ws.Range("B7").Value = "MaquinaTren"
$ws.Range("C7").Value = "cobrar"
$ws.Range("D7").Value = "descuenta el valor del boleto (Tren)"
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = "void"

# --- Las filas 8 a 11 (antes 7 a 10: Maquina.carga, Maquina.devolverSaldo,
#     SubeVirtual.calcularDescuento, BoletoTren.calcularValor) ya quedaron
#     correctamente desplazadas por la insercion de fila, sin cambios de contenido ---

# --- Nueva fila 12: BoletoColectivo.calcularValor ---
$ws.Range("B12").Value = "BoletoColectivo"
$ws.Range("C12").Value = "calcularValor"
$ws.Range("D12").Value = "calcula el valor del boleto según la seccion "
$ws.Range("E12").Value = "_"
$ws.Range("F12").Value = "float"

# --- Ajustar la celda seleccionada como en el archivo final ---
$ws.Range("E8").Select() | Out-Null
